$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Content.Find.Execute("2026-02-05 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-02-06 Friday", 2)

# Update the multiplication problems in the table.
# The table has 5 "data" rows (1, 5, 10, 15, 20) x 5 columns. Some values
# repeat (e.g. "239x3=") so addressing cells directly by row/column avoids
# ambiguity that a plain text Find/Replace would run into.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text  = "924×9="
$t.Cell(1, 2).Range.Text  = "425×5="
$t.Cell(1, 3).Range.Text  = "177×3="
$t.Cell(1, 4).Range.Text  = "549×8="
$t.Cell(1, 5).Range.Text  = "996×6="

$t.Cell(5, 1).Range.Text  = "647×6="
$t.Cell(5, 2).Range.Text  = "544×7="
$t.Cell(5, 3).Range.Text  = "774×9="
$t.Cell(5, 4).Range.Text  = "189×6="
$t.Cell(5, 5).Range.Text  = "922×2="

$t.Cell(10, 1).Range.Text = "796×4="
$t.Cell(10, 2).Range.Text = "873×5="
$t.Cell(10, 3).Range.Text = "780×4="
$t.Cell(10, 4).Range.Text = "661×3="
$t.Cell(10, 5).Range.Text = "803×9="

$t.Cell(15, 1).Range.Text = "555×5="
$t.Cell(15, 2).Range.Text = "722×4="
$t.Cell(15, 3).Range.Text = "483×9="
$t.Cell(15, 4).Range.Text = "797×3="
$t.Cell(15, 5).Range.Text = "371×9="

$t.Cell(20, 1).Range.Text = "753×7="
$t.Cell(20, 2).Range.Text = "778×2="
$t.Cell(20, 3).Range.Text = "379×3="
$t.Cell(20, 4).Range.Text = "182×2="
$t.Cell(20, 5).Range.Text = "337×8="
